$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Also added an input to allow for a change in trackwidth",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Also added an input to allow for a change in trackwidth. Changing the trackwidth allows for a more narrow extraction width to hopefully help get rid of any background interference with the TRUVOT spectra.",
    2)
